$d = $word.ActiveDocument

# 1. Créditos-aula: 4 -> 2
$d.Content.Find.Execute("Créditos-aula: 4", $true, $false, $false, $false, $false, $true, 1, $false, "Créditos-aula: 2", 2)

# 2. Carga horária: 60 h -> 30 h
$d.Content.Find.Execute("Carga horária: 60 h", $true, $false, $false, $false, $false, $true, 1, $false, "Carga horária: 30 h", 2)

# 3. Ativação: 01/01/2020 -> 01/01/2025
$d.Content.Find.Execute("Ativação: 01/01/2020", $true, $false, $false, $false, $false, $true, 1, $false, "Ativação: 01/01/2025", 2)

# 4. Shorten "Programa resumido" sentence
$d.Content.Find.Execute("Técnicas de Materialografia. Calorimetria e análises térmicas de materiais.", $true, $false, $false, $false, $false, $true, 1, $false, "Técnicas de Materialografia.", 2)

# 5. Remove the "2. CALORIMETRIA..." paragraph content (including the line break before it)
$oldProgramBlock = "Microscópio óptico de reflexão." + [char]11 + "2. CALORIMETRIA E ANÁLISES TÉRMICAS: Fundamentos termodinâmicos da calorimetria e análises térmicas. Princípios de calorimetria e tipos de calorímetros. Análise térmica diferencial (DTA) e calorimetria exploratória diferencial (DSC): princípios de DTA e DSC; tipos de equipamentos: DSC de compensação de energia e DSC de fluxo de calor. Aplicações de DTA e DSC. Equipamento: cadinhos de DTA; cadinhos de DSC. Cálculo de entalpia; linha base e cálculo de calor específico. Determinação de transição de fases. Determinação do diagrama de fases de ligas binárias por DTA/DSC. Cálculos cinéticos de cristalização, transições de fases e reações de polimerização. Termogravimetria (TGA): definição; aplicações da TGA. Equipamento: forno; programador de temperatura; termopar; balança; tipos de cadinho. Avaliação de estabilidade térmica e estudos de envelhecimento de polímeros. Técnicas de análises térmicas acopladas a análise de gases evolvidos por espectrometria de massa (TGA-MS) e FTIR (TGA-FTIR)."
$d.Content.Find.Execute($oldProgramBlock, $true, $false, $false, $false, $false, $true, 1, $false, "Microscópio óptico de reflexão.", 2)

# 6. Bibliography edits
$d.Content.Find.Execute("São Paulo – 1974", $true, $false, $false, $false, $false, $true, 1, $false, "SãoPaulo – 1974", 2)

$oldBib = "2008.AZEVEDO, A. D.; MOTHE, C. G. Análaise Térmica de Materiais. São Paulo: ARTLIBER, 2009.BROWN, M.E. Handbook of Thermal Analysis and Calorimetry, Amsterdam: Elsevier Science, 1998.HATAKEYAMA, T.; ZHENHAI, L. Handbook of Thermal Analysis, New York: Wiley, 1999.HAINES, P. J. Principles of Thermal Analysis and Calorimetry, Royal Society of Chemistry, 2002.MULLER, A. Solidificação e Análise Térmica dos Metais. Porto Alegre: Ed. UFRGS, 2002.SPEYER, R. Thermal analysis of materials, New York: Marcel Dekker, 1994.REED-HILL"
$d.Content.Find.Execute($oldBib, $true, $false, $false, $false, $false, $true, 1, $false, "2008.REED-HILL", 2)

$d.Content.Find.Execute("1982. Nondestructive", $true, $false, $false, $false, $false, $true, 1, $false, "1982.Nondestructive", 2)

$d.Content.Find.Execute("New York. YACOBI", $true, $false, $false, $false, $false, $true, 1, $false, "New York.YACOBI", 2)
